# Applies the scheduled-runner data refresh to the Leve profit tables.
# Columns H-N hold externally-sourced market data (currentAveragePrice*,
# LevePrice*, LeveProfit*) that the runner snapshots on each run; this
# script rewrites the affected cells per-sheet to their latest values.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 224.44444
$ws.Range("I2").Value = 241.25
$ws.Range("J2").Value = 90
$ws.Range("K2").Value = 241.25
$ws.Range("L2").Value = 90
$ws.Range("M2").Value = -128.25
$ws.Range("N2").Value = -316
$ws.Range("H62").Value = 1995
$ws.Range("J62").Value = 1995
$ws.Range("L62").Value = 1995
$ws.Range("N62").Value = -3243
$ws.Range("H64").Value = 9409.799999999999
$ws.Range("I64").Value = 8049.5
$ws.Range("J64").Value = 9749.875
$ws.Range("K64").Value = 8049.5
$ws.Range("L64").Value = 9749.875
$ws.Range("M64").Value = -7801.5
$ws.Range("N64").Value = -10245.875
$ws.Range("H65").Value = 1995
$ws.Range("J65").Value = 1995
$ws.Range("L65").Value = 9975
$ws.Range("N65").Value = -16215
$ws.Range("H67").Value = 9409.799999999999
$ws.Range("I67").Value = 8049.5
$ws.Range("J67").Value = 9749.875
$ws.Range("K67").Value = 8049.5
$ws.Range("L67").Value = 9749.875
$ws.Range("M67").Value = -7191.5
$ws.Range("N67").Value = -11465.875
$ws.Range("H70").Value = 3456.3635
$ws.Range("I70").Value = 3999
$ws.Range("J70").Value = 3402.1
$ws.Range("K70").Value = 11997
$ws.Range("L70").Value = 10206.3
$ws.Range("M70").Value = -11727
$ws.Range("N70").Value = -10746.3
$ws.Range("H73").Value = 3456.3635
$ws.Range("I73").Value = 3999
$ws.Range("J73").Value = 3402.1
$ws.Range("K73").Value = 11997
$ws.Range("L73").Value = 10206.3
$ws.Range("M73").Value = -11061
$ws.Range("N73").Value = -12078.3
$ws.Range("H74").Value = 12166.667
$ws.Range("I74").Value = 11500
$ws.Range("J74").Value = 12500
$ws.Range("K74").Value = 11500
$ws.Range("L74").Value = 12500
$ws.Range("M74").Value = -10564
$ws.Range("N74").Value = -14372
$ws.Range("H77").Value = 12166.667
$ws.Range("I77").Value = 11500
$ws.Range("J77").Value = 12500
$ws.Range("K77").Value = 57500
$ws.Range("L77").Value = 62500
$ws.Range("M77").Value = -52820
$ws.Range("N77").Value = -71860
$ws.Range("H87").Value = 26584.5
$ws.Range("J87").Value = 26584.5
$ws.Range("L87").Value = 26584.5
$ws.Range("N87").Value = -29080.5
$ws.Range("H90").Value = 26584.5
$ws.Range("J90").Value = 26584.5
$ws.Range("L90").Value = 79753.5
$ws.Range("N90").Value = -92233.5
$ws.Range("H123").Value = 59187.445
$ws.Range("J123").Value = 59999.59
$ws.Range("L123").Value = 59999.59
$ws.Range("N123").Value = -69799.59
$ws.Range("H141").Value = 4436.353
$ws.Range("I141").Value = 3151.125
$ws.Range("K141").Value = 9453.375
$ws.Range("M141").Value = -4273.375

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 495
$ws.Range("I3").Value = 495
$ws.Range("K3").Value = 495
$ws.Range("M3").Value = -380
$ws.Range("H32").Value = 52643730
$ws.Range("I32").Value = 52643730
$ws.Range("K32").Value = 52643730
$ws.Range("M32").Value = -52643443

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6832.3335
$ws.Range("I20").Value = 8999
$ws.Range("K20").Value = 8999
$ws.Range("M20").Value = -8752
$ws.Range("H86").Value = 5372.364
$ws.Range("I86").Value = 4682.8335
$ws.Range("K86").Value = 4682.8335
$ws.Range("M86").Value = -3559.8335
$ws.Range("H89").Value = 5372.364
$ws.Range("I89").Value = 4682.8335
$ws.Range("K89").Value = 23414.1675
$ws.Range("M89").Value = -17798.1675
$ws.Range("H99").Value = 1791.1351
$ws.Range("I99").Value = 1723.48
$ws.Range("J99").Value = 1932.0834
$ws.Range("K99").Value = 1723.48
$ws.Range("L99").Value = 1932.0834
$ws.Range("M99").Value = -225.48
$ws.Range("N99").Value = -4928.0834
$ws.Range("H107").Value = 1993
$ws.Range("I107").Value = 1993
$ws.Range("K107").Value = 1993
$ws.Range("M107").Value = -73

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4008.9
$ws.Range("I62").Value = 4010
$ws.Range("J62").Value = 3999
$ws.Range("K62").Value = 4010
$ws.Range("L62").Value = 3999
$ws.Range("M62").Value = -3386
$ws.Range("N62").Value = -5247
$ws.Range("H65").Value = 4008.9
$ws.Range("I65").Value = 4010
$ws.Range("J65").Value = 3999
$ws.Range("K65").Value = 20050
$ws.Range("L65").Value = 3999
$ws.Range("M65").Value = -16930
$ws.Range("N65").Value = -26235

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8666.608
$ws.Range("J70").Value = 10447.454
$ws.Range("L70").Value = 10447.454
$ws.Range("N70").Value = -10987.454
$ws.Range("H73").Value = 8666.608
$ws.Range("J73").Value = 10447.454
$ws.Range("L73").Value = 10447.454
$ws.Range("N73").Value = -12319.454
$ws.Range("H132").Value = 930.3333
$ws.Range("I132").Value = 930.3333
$ws.Range("K132").Value = 2790.9999
$ws.Range("M132").Value = -260.9998999999998

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4480.75
$ws.Range("I61").Value = 4316
$ws.Range("J61").Value = 4579.6
$ws.Range("K61").Value = 4316
$ws.Range("L61").Value = 4579.6
$ws.Range("M61").Value = -4114
$ws.Range("N61").Value = -4983.6
$ws.Range("H82").Value = 4815.3335
$ws.Range("I82").Value = 2592.5
$ws.Range("J82").Value = 7355.7144
$ws.Range("K82").Value = 2592.5
$ws.Range("L82").Value = 7355.7144
$ws.Range("M82").Value = -2231.5
$ws.Range("N82").Value = -8077.7144
$ws.Range("H85").Value = 4815.3335
$ws.Range("I85").Value = 2592.5
$ws.Range("J85").Value = 7355.7144
$ws.Range("K85").Value = 2592.5
$ws.Range("L85").Value = 7355.7144
$ws.Range("M85").Value = -1344.5
$ws.Range("N85").Value = -9851.714400000001
$ws.Range("H100").Value = 4100.65
$ws.Range("I100").Value = 2078.7778
$ws.Range("K100").Value = 2078.7778
$ws.Range("M100").Value = -1537.7778
$ws.Range("H113").Value = 4480.75
$ws.Range("I113").Value = 4316
$ws.Range("J113").Value = 4579.6
$ws.Range("K113").Value = 4316
$ws.Range("L113").Value = 4579.6
$ws.Range("M113").Value = -2146
$ws.Range("N113").Value = -8919.6

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H9").Value = 7500
$ws.Range("I9").Value = 2000
$ws.Range("K9").Value = 2000
$ws.Range("M9").Value = -1860
$ws.Range("H14").Value = 5503.25
$ws.Range("I14").Value = 5503.25
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 5503.25
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -5335.25
$ws.Range("N14").ClearContents()
$ws.Range("H107").Value = 1500
$ws.Range("I107").Value = 1500
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 4500
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -2580
$ws.Range("N107").ClearContents()
$ws.Range("H132").Value = 2118.2068
$ws.Range("I132").Value = 1074.0526
$ws.Range("K132").Value = 3222.1578
$ws.Range("M132").Value = -692.1578
